$wb = $excel.ActiveWorkbook

# --- Update Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-05-13T18:56:15+00:00"
$meta.Range("B22").NumberFormat = "@"
$meta.Range("B22").Value = "15"
# Restore original cell formatting (the NumberFormat tweak above is only
# needed transiently so Excel stores the value as text instead of a number).
$meta.Range("B21").Copy()
$meta.Range("B22").PasteSpecial(-4122)

# --- Append new concept rows to Concepts sheet ---
$concepts = $wb.Worksheets.Item("Concepts")

$newRows = @(
    @("1", "RNAS", "RNA-Seq"),
    @("1", "CHIPS", "ChIP-Seq"),
    @("1", "ATACS", "ATAC-Seq"),
    @("1", "HIC", "Hi-C"),
    @("1", "BIS", "Bisulfite-Seq"),
    @("1", "TMS", "Targeted Methyl-Seq"),
    @("1", "MRES", "MRE-Seq"),
    @("1", "MDIPS", "MeDIP-Seq"),
    @("1", "MCCS", "MCC-Seq"),
    @("1", "MBDS", "MBD-Seq"),
    @("1", "EMS", "EM-Seq"),
    @("1", "DMS", "Direct Methyl-Seq")
)

$startRow = 5

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Force the Level column to be stored as text ("1") rather than a
    # number, matching the existing rows above it.
    $concepts.Cells.Item($r, 1).NumberFormat = "@"
    $concepts.Cells.Item($r, 1).Value = $row[0]
    $concepts.Cells.Item($r, 2).Value = $row[1]
    $concepts.Cells.Item($r, 3).Value = $row[2]
    $concepts.Cells.Item($r, 4).Value = ""

    # Copy the formatting from the row above (style index matches the
    # rest of the table) without touching the values just entered.
    $concepts.Range("A4:D4").Copy()
    $concepts.Range("A$r`:D$r").PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
